# Generate Report for Archive
# Replace the "Ready for handoff" localization status with "In Translation"
# across the Overview summary sheet and each per-language detail sheet, then
# re-tighten the affected "Status" columns so they fit the (shorter) new text.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

# --- Overview sheet: per-language status columns (E = zh-cn, F = de-de) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus
$wsOverview.Range("E4").Value = $newStatus
$wsOverview.Range("F4").Value = $newStatus

# --- Per-language detail sheets: "Status" column (column C) ---
$languageSheets = @("zh-cn", "de-de")
foreach ($sheetName in $languageSheets) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("C2").Value = $newStatus
    $ws.Range("C3").Value = $newStatus
    $ws.Range("C4").Value = $newStatus
}

# --- Re-fit the columns that held the old/new status text ---
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
foreach ($sheetName in $languageSheets) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Columns.Item(3).ColumnWidth = 12.5
}
